$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SV4's student ID (row 13 -> D13) and SV3's student ID (row 12 -> D12)
# Written in this order so the shared-string table gains "PS09095" before
# "PS09117", matching the target workbook's string indices.
$ws.Range("D13").Value = "PS09095"
$ws.Range("D12").Value = "PS09117"

# Matches the final selection state recorded in the saved sheet view
$null = $ws.Range("D13").Select()
